# Update "想去人数" (interested-attendee count) values by +1 for five events
# that appear both on the "展览" (Exhibition) sheet and the aggregated
# "全部类型" (All Types) sheet, matching the upstream data refresh.

$wb = $excel.ActiveWorkbook

# Sheet "展览": rows 5, 9, 19, 24, 29 -> column F
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value  = 1042
$wsExhibit.Range("F9").Value  = 1488
$wsExhibit.Range("F19").Value = 273
$wsExhibit.Range("F24").Value = 39
$wsExhibit.Range("F29").Value = 47

# Sheet "全部类型": rows 16, 20, 30, 37, 44 -> column F (same events as above)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F16").Value = 1042
$wsAll.Range("F20").Value = 1488
$wsAll.Range("F30").Value = 273
$wsAll.Range("F37").Value = 39
$wsAll.Range("F44").Value = 47
